# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 with the latest scraped figures.
#
# Column D prices are stored as plain text (thousands-separated with
# dots, e.g. "43.186.78"), so values that happen to look like plain
# numbers (e.g. "114.08") must be forced to text -- otherwise Excel
# auto-converts them to numeric cells. We do this by switching the
# cell's NumberFormat to Text ("@") before writing, then restoring the
# cell's style back to "Normal" afterwards so no visible formatting
# change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) value or $null, whether D needs to be
# forced to text, new Volume(1h) (E) value or $null.
$updates = @(
    @(2, '43.186.78', $false, '  +1.32%  '),
    @(3, '2.280.83', $false, '  +1.91%  '),
    @(4, $null, $false, '  +0.28%  '),
    @(5, '114.08', $true, '  -0.61%  '),
    @(6, '304.57', $true, '  +6.46%  '),
    @(7, '0.631', $true, '  +0.79%  '),
    @(8, $null, $false, '  -0.12%  '),
    @(9, '0.615', $true, '  +0.32%  '),
    @(10, '44.74', $true, '  -4.17%  '),
    @(11, '0.0927', $true, '  -0.42%  '),
    @(12, $null, $false, '  +1.15%  '),
    @(13, $null, $false, '  -2.74%  '),
    @(14, '1.06', $true, '  +19.76%  '),
    @(15, $null, $false, '  +0.06%  '),
    @(16, $null, $false, '  +0.38%  '),
    @(17, '2.620.75', $false, '  +1.73%  '),
    @(18, '2.282.56', $false, '  +2.04%  '),
    @(19, '43.109.18', $false, '  +0.97%  '),
    @(20, $null, $false, '  +0.12%  '),
    @(21, $null, $false, '  +4.60%  '),
    @(22, '75.05', $true, '  +2.28%  '),
    @(23, $null, $false, '  +11.88%  '),
    @(24, '2.46', $true, '  +4.22%  '),
    @(25, '254.79', $true, '  +9.77%  '),
    @(26, $null, $false, '  -1.83%  '),
    @(27, $null, $false, '  -3.40%  '),
    @(28, '0.999', $true, '  -0.33%  '),
    @(29, $null, $false, '  +0.33%  '),
    @(30, '38.21', $true, '  -4.69%  '),
    @(31, '175.25', $true, '  -0.12%  '),
    @(32, '22.20', $true, '  +4.79%  '),
    @(33, $null, $false, '  -3.80%  '),
    @(34, '0.0899', $true, '  -0.51%  '),
    @(35, '5.70', $true, '  +2.05%  '),
    @(36, '5.05', $true, '  +9.30%  '),
    @(37, $null, $false, '  +0.76%  '),
    @(38, '4.25', $true, '  -8.03%  '),
    @(39, '0.0379', $true, '  +1.98%  '),
    @(40, $null, $false, '  -0.64%  '),
    @(41, '2.54', $true, '  -3.32%  '),
    @(42, '72.70', $true, '  -0.08%  '),
    @(43, $null, $false, '  -0.40%  '),
    @(44, $null, $false, '  +0.22%  '),
    @(45, '12.66', $true, '  -6.20%  '),
    @(46, '1.38', $true, '  +3.86%  '),
    @(47, '5.64', $true, '  +0.63%  '),
    @(48, '106.76', $true, '  +5.54%  '),
    @(49, '8.80', $true, '  +2.90%  '),
    @(50, '1.29', $true, '  -0.17%  '),
    @(51, '74.22', $true, '  +6.51%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $priceVal = $u[1]
    $forceText = $u[2]
    $volVal = $u[3]

    if ($null -ne $priceVal) {
        $cell = $ws.Cells.Item($row, 4)
        if ($forceText) {
            $cell.NumberFormat = "@"
            $cell.Value = $priceVal
            $cell.Style = "Normal"
        } else {
            $cell.Value = $priceVal
        }
    }
    if ($null -ne $volVal) {
        $ws.Cells.Item($row, 5).Value = $volVal
    }
}
